# Gruppendaten.xlsx - "Anpassung für Excel Export der Gruppen"
# Rebuild the sheet: new title, summary/KPI block, and a real Excel Table
# ("Tabelle1") with per-column SUM() formulas feeding the KPI cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe the old header row content (row 2) - it gets replaced further down.
# ---------------------------------------------------------------------------
$ws.Range("A2:K2").Clear()

# ---------------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Verbindliche Teilnehmende für das Zeltlager"

# ---------------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 39.85546875
$ws.Columns.Item(2).ColumnWidth = 18.7109375
$ws.Columns.Item(3).ColumnWidth = 18.7109375
$ws.Columns.Item(4).ColumnWidth = 16
$ws.Columns.Item(5).ColumnWidth = 22.140625
$ws.Columns.Item(6).ColumnWidth = 31.85546875
$ws.Columns.Item(7).ColumnWidth = 20.5703125
$ws.Columns.Item(8).ColumnWidth = 45.140625
$ws.Columns.Item(9).ColumnWidth = 18.85546875
$ws.Columns.Item(10).ColumnWidth = 31.7109375
$ws.Columns.Item(11).ColumnWidth = 29.85546875
$ws.Columns.Item(12).ColumnWidth = 28.5703125
$ws.Columns.Item(13).ColumnWidth = 35.42578125
$ws.Columns.Item(14).ColumnWidth = 22.5703125
$ws.Columns.Item(15).ColumnWidth = 28.140625

# ---------------------------------------------------------------------------
# 4. KPI / summary block (rows 2-5)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 15.75

$ws.Range("A3").Value = "Anzahl an Teilnehmdenen:"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Size = 12
$ws.Range("A3").HorizontalAlignment = -4152
$ws.Range("B3").Value = 0
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Font.Size = 12
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("A3:B3").Borders.Item(8).LineStyle = 1
$ws.Range("A3:B3").Borders.Item(8).Weight = 4
$ws.Range("A3:B3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:B3").Borders.Item(9).Weight = 4
$ws.Range("A3").Borders.Item(7).LineStyle = 1
$ws.Range("A3").Borders.Item(7).Weight = 4
$ws.Range("B3").Borders.Item(10).LineStyle = 1
$ws.Range("B3").Borders.Item(10).Weight = 4
$ws.Rows.Item(3).RowHeight = 16.5
$ws.Range("C3").HorizontalAlignment = -4131

$ws.Rows.Item(4).RowHeight = 15.75

$ws.Range("A5").Value = "Davon mit besonderen Essgewohnheiten:"
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").HorizontalAlignment = -4152
$ws.Range("B5").Value = 0
$ws.Range("B5").Font.Bold = $true
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("A5:B5").Borders.Item(8).LineStyle = 1
$ws.Range("A5:B5").Borders.Item(8).Weight = 4
$ws.Range("A5:B5").Borders.Item(9).LineStyle = 1
$ws.Range("A5:B5").Borders.Item(9).Weight = 4
$ws.Range("A5").Borders.Item(7).LineStyle = 1
$ws.Range("A5").Borders.Item(7).Weight = 4
$ws.Range("B5").Borders.Item(10).LineStyle = 1
$ws.Range("B5").Borders.Item(10).Weight = 4
$ws.Range("C5").HorizontalAlignment = -4131

$ws.Range("E5").Value = "Davon Vegetrarier:"
$ws.Range("E5").HorizontalAlignment = -4152
$ws.Range("F5").Value = 0
$ws.Range("F5").HorizontalAlignment = -4131
$ws.Range("E5").Borders.Item(7).LineStyle = 1
$ws.Range("E5").Borders.Item(7).Weight = 4
$ws.Range("E5").Borders.Item(8).LineStyle = 1
$ws.Range("E5").Borders.Item(8).Weight = 4
$ws.Range("E5").Borders.Item(9).LineStyle = 1
$ws.Range("E5").Borders.Item(9).Weight = 4
$ws.Range("F5").Borders.Item(10).LineStyle = 1
$ws.Range("F5").Borders.Item(10).Weight = 4
$ws.Range("F5").Borders.Item(8).LineStyle = 1
$ws.Range("F5").Borders.Item(8).Weight = 4
$ws.Range("F5").Borders.Item(9).LineStyle = 1
$ws.Range("F5").Borders.Item(9).Weight = 4

$ws.Range("H5").Value = "Davon Veganer:"
$ws.Range("H5").HorizontalAlignment = -4152
$ws.Range("I5").Value = 0
$ws.Range("I5").HorizontalAlignment = -4131
$ws.Range("H5").Borders.Item(7).LineStyle = 1
$ws.Range("H5").Borders.Item(7).Weight = 4
$ws.Range("H5").Borders.Item(8).LineStyle = 1
$ws.Range("H5").Borders.Item(8).Weight = 4
$ws.Range("H5").Borders.Item(9).LineStyle = 1
$ws.Range("H5").Borders.Item(9).Weight = 4
$ws.Range("I5").Borders.Item(10).LineStyle = 1
$ws.Range("I5").Borders.Item(10).Weight = 4
$ws.Range("I5").Borders.Item(8).LineStyle = 1
$ws.Range("I5").Borders.Item(8).Weight = 4
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 4

# ---------------------------------------------------------------------------
# 5. Table header row (row 10) + the actual ListObject/Table "Tabelle1"
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Name"
$ws.Range("B10").Value = "Vorname"
$ws.Range("C10").Value = "Geburtsdatum"
$ws.Range("D10").Value = "Alter"
$ws.Range("E10").Value = "M/W/D/N"
$ws.Range("F10").Value = "PLZ"
$ws.Range("G10").Value = "Ort"
$ws.Range("H10").Value = "Straße, Hausnummer"
$ws.Range("I10").Value = "TN Status"
$ws.Range("J10").Value = "Jugendfeuerwehr"
$ws.Range("K10").Value = "Organisationseinheit"
$ws.Range("L10").Value = "Essgewohnheiten"
$ws.Range("M10").Value = "Unverträglichkeiten"
$ws.Range("N10").Value = "Zu überweisen"
$ws.Range("O10").Value = "Bereits überwiesen"

$lo = $ws.ListObjects.Add(1, $ws.Range("A10:O11"), $null, 1)
$lo.Name = "Tabelle1"
$lo.TableStyle = "TableStyleMedium2"

# ---------------------------------------------------------------------------
# 6. Row 9 - per-column SUM() totals driven off the table columns
# ---------------------------------------------------------------------------
$ws.Range("A9").Formula = "=SUM(Tabelle1[Name])"
$ws.Range("B9").Formula = "=SUM(Tabelle1[Vorname])"
$ws.Range("C9").Formula = "=SUM(Tabelle1[Geburtsdatum])"
$ws.Range("D9").Formula = "=SUM(Tabelle1[Alter])"
$ws.Range("E9").Formula = "=SUM(Tabelle1[M/W/D/N])"
$ws.Range("F9").Formula = "=SUM(Tabelle1[PLZ])"
$ws.Range("G9").Formula = "=SUM(Tabelle1[Ort])"
$ws.Range("H9").Formula = "=SUM(Tabelle1[Straße, Hausnummer])"
$ws.Range("I9").Formula = "=SUM(Tabelle1[TN Status])"
$ws.Range("J9").Formula = "=SUM(Tabelle1[Jugendfeuerwehr])"
$ws.Range("K9").Formula = "=SUM(Tabelle1[Organisationseinheit])"
$ws.Range("L9").Formula = "=SUM(Tabelle1[Essgewohnheiten])"
$ws.Range("M9").Formula = "=SUM(Tabelle1[Unverträglichkeiten])"
$ws.Range("N9").Formula = "=SUM(Tabelle1[Zu überweisen])"
$ws.Range("O9").Formula = "=SUM(Tabelle1[Bereits überwiesen])"
$ws.Range("A9:O9").Font.Bold = $true

# ---------------------------------------------------------------------------
# 7. Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("A9").Select()
